# edit.ps1
# Applies the "optimization_results" update:
#  - Results sheet: re-optimized award table, rows 2-15 updated and two new
#    rows (16, 17) added for the split Facility4/Facility5 bids (dimension
#    A1:O15 -> A1:O17).
#  - LP Model sheet: Rule_0_* constraints rewritten from per-bid capacity
#    caps ("x_X_n <= 300") to minimum-award-split rules
#    ("x_A_n + x_B_n >= 13" etc), matching the "extended rule" text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Results sheet - rewrite the data rows (A2:O17)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Results")

    # Row 2
    $ws.Cells.Item(2,1).Value = 1
    $ws.Cells.Item(2,2).Value = "A"
    $ws.Cells.Item(2,3).Value = "Facility1"
    $ws.Cells.Item(2,4).Value = "A"
    $ws.Cells.Item(2,5).Value = 100
    $ws.Cells.Item(2,6).Value = 68700
    $ws.Cells.Item(2,7).Value = "A"
    $ws.Cells.Item(2,8).Value = 20
    $c = $ws.Cells.Item(2,9)
    $c.NumberFormat = "@"
    $c.Value = "1%"
    $c.ClearFormats()
    $ws.Cells.Item(2,10).Value = 19.8
    $ws.Cells.Item(2,11).Value = 13602.6
    $ws.Cells.Item(2,12).Value = 687
    $ws.Cells.Item(2,13).Value = 55097.4
    $c = $ws.Cells.Item(2,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(2,15).Value = 0
    # Row 3
    $ws.Cells.Item(3,1).Value = 1
    $ws.Cells.Item(3,2).Value = "B"
    $ws.Cells.Item(3,3).Value = "Facility1"
    $ws.Cells.Item(3,4).Value = "A"
    $ws.Cells.Item(3,5).Value = 100
    $ws.Cells.Item(3,6).Value = 1300
    $ws.Cells.Item(3,7).Value = "C"
    $ws.Cells.Item(3,8).Value = 55
    $c = $ws.Cells.Item(3,9)
    $c.NumberFormat = "@"
    $c.Value = "4%"
    $c.ClearFormats()
    $ws.Cells.Item(3,10).Value = 52.8
    $ws.Cells.Item(3,11).Value = 686.4
    $ws.Cells.Item(3,12).Value = 13
    $ws.Cells.Item(3,13).Value = 613.6
    $c = $ws.Cells.Item(3,14)
    $c.NumberFormat = "@"
    $c.Value = "7%"
    $c.ClearFormats()
    $ws.Cells.Item(3,15).Value = 48.048
    # Row 4
    $ws.Cells.Item(4,1).Value = 2
    $ws.Cells.Item(4,2).Value = "A"
    $ws.Cells.Item(4,3).Value = "Facility1"
    $ws.Cells.Item(4,4).Value = "B"
    $ws.Cells.Item(4,5).Value = 156
    $ws.Cells.Item(4,6).Value = 1404000
    $ws.Cells.Item(4,7).Value = "C"
    $ws.Cells.Item(4,8).Value = 75
    $c = $ws.Cells.Item(4,9)
    $c.NumberFormat = "@"
    $c.Value = "4%"
    $c.ClearFormats()
    $ws.Cells.Item(4,10).Value = 72
    $ws.Cells.Item(4,11).Value = 648000
    $ws.Cells.Item(4,12).Value = 9000
    $ws.Cells.Item(4,13).Value = 756000
    $c = $ws.Cells.Item(4,14)
    $c.NumberFormat = "@"
    $c.Value = "7%"
    $c.ClearFormats()
    $ws.Cells.Item(4,15).Value = 45360.00000000001
    # Row 5
    $ws.Cells.Item(5,1).Value = 3
    $ws.Cells.Item(5,2).Value = "A"
    $ws.Cells.Item(5,3).Value = "Facility1"
    $ws.Cells.Item(5,4).Value = "C"
    $ws.Cells.Item(5,5).Value = 423
    $ws.Cells.Item(5,6).Value = 126900
    $ws.Cells.Item(5,7).Value = "A"
    $ws.Cells.Item(5,8).Value = 55
    $c = $ws.Cells.Item(5,9)
    $c.NumberFormat = "@"
    $c.Value = "1%"
    $c.ClearFormats()
    $ws.Cells.Item(5,10).Value = 54.45
    $ws.Cells.Item(5,11).Value = 16335
    $ws.Cells.Item(5,12).Value = 300
    $ws.Cells.Item(5,13).Value = 110565
    $c = $ws.Cells.Item(5,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(5,15).Value = 0
    # Row 6
    $ws.Cells.Item(6,1).Value = 3
    $ws.Cells.Item(6,2).Value = "B"
    $ws.Cells.Item(6,3).Value = "Facility1"
    $ws.Cells.Item(6,4).Value = "C"
    $ws.Cells.Item(6,5).Value = 423
    $ws.Cells.Item(6,6).Value = 126900
    $ws.Cells.Item(6,7).Value = "C"
    $ws.Cells.Item(6,8).Value = 60
    $c = $ws.Cells.Item(6,9)
    $c.NumberFormat = "@"
    $c.Value = "4%"
    $c.ClearFormats()
    $ws.Cells.Item(6,10).Value = 57.59999999999999
    $ws.Cells.Item(6,11).Value = 17280
    $ws.Cells.Item(6,12).Value = 300
    $ws.Cells.Item(6,13).Value = 109620
    $c = $ws.Cells.Item(6,14)
    $c.NumberFormat = "@"
    $c.Value = "7%"
    $c.ClearFormats()
    $ws.Cells.Item(6,15).Value = 1209.6
    # Row 7
    $ws.Cells.Item(7,1).Value = 4
    $ws.Cells.Item(7,2).Value = "A"
    $ws.Cells.Item(7,3).Value = "Facility2"
    $ws.Cells.Item(7,4).Value = "C"
    $ws.Cells.Item(7,5).Value = 453
    $ws.Cells.Item(7,6).Value = 2562621
    $ws.Cells.Item(7,7).Value = "C"
    $ws.Cells.Item(7,8).Value = 19
    $c = $ws.Cells.Item(7,9)
    $c.NumberFormat = "@"
    $c.Value = "4%"
    $c.ClearFormats()
    $ws.Cells.Item(7,10).Value = 18.24
    $ws.Cells.Item(7,11).Value = 103183.68
    $ws.Cells.Item(7,12).Value = 5657
    $ws.Cells.Item(7,13).Value = 2459437.32
    $c = $ws.Cells.Item(7,14)
    $c.NumberFormat = "@"
    $c.Value = "7%"
    $c.ClearFormats()
    $ws.Cells.Item(7,15).Value = 7222.8576
    # Row 8
    $ws.Cells.Item(8,1).Value = 4
    $ws.Cells.Item(8,2).Value = "B"
    $ws.Cells.Item(8,3).Value = "Facility2"
    $ws.Cells.Item(8,4).Value = "C"
    $ws.Cells.Item(8,5).Value = 453
    $ws.Cells.Item(8,6).Value = 5889
    $ws.Cells.Item(8,7).Value = "A"
    $ws.Cells.Item(8,8).Value = 23
    $c = $ws.Cells.Item(8,9)
    $c.NumberFormat = "@"
    $c.Value = "1%"
    $c.ClearFormats()
    $ws.Cells.Item(8,10).Value = 22.77
    $ws.Cells.Item(8,11).Value = 296.01
    $ws.Cells.Item(8,12).Value = 13
    $ws.Cells.Item(8,13).Value = 5592.99
    $c = $ws.Cells.Item(8,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(8,15).Value = 0
    # Row 9
    $ws.Cells.Item(9,1).Value = 5
    $ws.Cells.Item(9,2).Value = "A"
    $ws.Cells.Item(9,3).Value = "Facility2"
    $ws.Cells.Item(9,4).Value = "C"
    $ws.Cells.Item(9,5).Value = 342
    $ws.Cells.Item(9,6).Value = 10944
    $ws.Cells.Item(9,7).Value = "C"
    $ws.Cells.Item(9,8).Value = 24
    $c = $ws.Cells.Item(9,9)
    $c.NumberFormat = "@"
    $c.Value = "4%"
    $c.ClearFormats()
    $ws.Cells.Item(9,10).Value = 23.04
    $ws.Cells.Item(9,11).Value = 737.28
    $ws.Cells.Item(9,12).Value = 32
    $ws.Cells.Item(9,13).Value = 10206.72
    $c = $ws.Cells.Item(9,14)
    $c.NumberFormat = "@"
    $c.Value = "7%"
    $c.ClearFormats()
    $ws.Cells.Item(9,15).Value = 51.6096
    # Row 10
    $ws.Cells.Item(10,1).Value = 5
    $ws.Cells.Item(10,2).Value = "B"
    $ws.Cells.Item(10,3).Value = "Facility2"
    $ws.Cells.Item(10,4).Value = "C"
    $ws.Cells.Item(10,5).Value = 342
    $ws.Cells.Item(10,6).Value = 4446
    $ws.Cells.Item(10,7).Value = "B"
    $ws.Cells.Item(10,8).Value = 34
    $c = $ws.Cells.Item(10,9)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(10,10).Value = 34
    $ws.Cells.Item(10,11).Value = 442
    $ws.Cells.Item(10,12).Value = 13
    $ws.Cells.Item(10,13).Value = 4004
    $c = $ws.Cells.Item(10,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(10,15).Value = 0
    # Row 11
    $ws.Cells.Item(11,1).Value = 6
    $ws.Cells.Item(11,2).Value = "A"
    $ws.Cells.Item(11,3).Value = "Facility2"
    $ws.Cells.Item(11,4).Value = "C"
    $ws.Cells.Item(11,5).Value = 653
    $ws.Cells.Item(11,6).Value = 144966
    $ws.Cells.Item(11,7).Value = "B"
    $ws.Cells.Item(11,8).Value = 24
    $c = $ws.Cells.Item(11,9)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(11,10).Value = 24
    $ws.Cells.Item(11,11).Value = 5328
    $ws.Cells.Item(11,12).Value = 222
    $ws.Cells.Item(11,13).Value = 139638
    $c = $ws.Cells.Item(11,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(11,15).Value = 0
    # Row 12
    $ws.Cells.Item(12,1).Value = 6
    $ws.Cells.Item(12,2).Value = "B"
    $ws.Cells.Item(12,3).Value = "Facility2"
    $ws.Cells.Item(12,4).Value = "C"
    $ws.Cells.Item(12,5).Value = 653
    $ws.Cells.Item(12,6).Value = 13060
    $ws.Cells.Item(12,7).Value = "A"
    $ws.Cells.Item(12,8).Value = 42
    $c = $ws.Cells.Item(12,9)
    $c.NumberFormat = "@"
    $c.Value = "1%"
    $c.ClearFormats()
    $ws.Cells.Item(12,10).Value = 41.58
    $ws.Cells.Item(12,11).Value = 831.5999999999999
    $ws.Cells.Item(12,12).Value = 20
    $ws.Cells.Item(12,13).Value = 12228.4
    $c = $ws.Cells.Item(12,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(12,15).Value = 0
    # Row 13
    $ws.Cells.Item(13,1).Value = 7
    $ws.Cells.Item(13,2).Value = "A"
    $ws.Cells.Item(13,3).Value = "Facility2"
    $ws.Cells.Item(13,4).Value = "C"
    $ws.Cells.Item(13,5).Value = 432
    $ws.Cells.Item(13,6).Value = 286848
    $ws.Cells.Item(13,7).Value = "A"
    $ws.Cells.Item(13,8).Value = 23
    $c = $ws.Cells.Item(13,9)
    $c.NumberFormat = "@"
    $c.Value = "1%"
    $c.ClearFormats()
    $ws.Cells.Item(13,10).Value = 22.77
    $ws.Cells.Item(13,11).Value = 15119.28
    $ws.Cells.Item(13,12).Value = 664
    $ws.Cells.Item(13,13).Value = 271728.72
    $c = $ws.Cells.Item(13,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(13,15).Value = 0
    # Row 14
    $ws.Cells.Item(14,1).Value = 8
    $ws.Cells.Item(14,2).Value = "A"
    $ws.Cells.Item(14,3).Value = "Facility3"
    $ws.Cells.Item(14,4).Value = "C"
    $ws.Cells.Item(14,5).Value = 456
    $ws.Cells.Item(14,6).Value = 10944
    $ws.Cells.Item(14,7).Value = "B"
    $ws.Cells.Item(14,8).Value = 13
    $c = $ws.Cells.Item(14,9)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(14,10).Value = 13
    $ws.Cells.Item(14,11).Value = 312
    $ws.Cells.Item(14,12).Value = 24
    $ws.Cells.Item(14,13).Value = 10632
    $c = $ws.Cells.Item(14,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(14,15).Value = 0
    # Row 15
    $ws.Cells.Item(15,1).Value = 9
    $ws.Cells.Item(15,2).Value = "A"
    $ws.Cells.Item(15,3).Value = "Facility4"
    $ws.Cells.Item(15,4).Value = "C"
    $ws.Cells.Item(15,5).Value = 234
    $ws.Cells.Item(15,6).Value = 51246
    $ws.Cells.Item(15,7).Value = "C"
    $ws.Cells.Item(15,8).Value = 13
    $c = $ws.Cells.Item(15,9)
    $c.NumberFormat = "@"
    $c.Value = "4%"
    $c.ClearFormats()
    $ws.Cells.Item(15,10).Value = 12.48
    $ws.Cells.Item(15,11).Value = 2733.12
    $ws.Cells.Item(15,12).Value = 219
    $ws.Cells.Item(15,13).Value = 48512.88
    $c = $ws.Cells.Item(15,14)
    $c.NumberFormat = "@"
    $c.Value = "7%"
    $c.ClearFormats()
    $ws.Cells.Item(15,15).Value = 191.3184
    # Row 16
    $ws.Cells.Item(16,1).Value = 9
    $ws.Cells.Item(16,2).Value = "B"
    $ws.Cells.Item(16,3).Value = "Facility4"
    $ws.Cells.Item(16,4).Value = "C"
    $ws.Cells.Item(16,5).Value = 234
    $ws.Cells.Item(16,6).Value = 3042
    $ws.Cells.Item(16,7).Value = "B"
    $ws.Cells.Item(16,8).Value = 56
    $c = $ws.Cells.Item(16,9)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(16,10).Value = 56
    $ws.Cells.Item(16,11).Value = 728
    $ws.Cells.Item(16,12).Value = 13
    $ws.Cells.Item(16,13).Value = 2314
    $c = $ws.Cells.Item(16,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(16,15).Value = 0
    # Row 17
    $ws.Cells.Item(17,1).Value = 10
    $ws.Cells.Item(17,2).Value = "A"
    $ws.Cells.Item(17,3).Value = "Facility5"
    $ws.Cells.Item(17,4).Value = "C"
    $ws.Cells.Item(17,5).Value = 231
    $ws.Cells.Item(17,6).Value = 3003
    $ws.Cells.Item(17,7).Value = "B"
    $ws.Cells.Item(17,8).Value = 13
    $c = $ws.Cells.Item(17,9)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(17,10).Value = 13
    $ws.Cells.Item(17,11).Value = 169
    $ws.Cells.Item(17,12).Value = 13
    $ws.Cells.Item(17,13).Value = 2834
    $c = $ws.Cells.Item(17,14)
    $c.NumberFormat = "@"
    $c.Value = "0%"
    $c.ClearFormats()
    $ws.Cells.Item(17,15).Value = 0

# ---------------------------------------------------------------------
# 2. LP Model sheet - update the extended Rule_0_* constraint text
# ---------------------------------------------------------------------
$lpWs = $wb.Worksheets.Item("LP Model")
$lpCell = $lpWs.Range("A2")
$lpText = $lpCell.Value()

    $lpText = $lpText.Replace('Rule_0_1: x_A_1 <= 300', 'Rule_0_1: x_B_1 + x_C_1 >= 13')
    $lpText = $lpText.Replace('Rule_0_10: x_C_10 <= 300', 'Rule_0_10: x_A_10 + x_B_10 >= 13')
    $lpText = $lpText.Replace('Rule_0_2: x_B_2 <= 300', 'Rule_0_2: x_A_2 + x_C_2 >= 13')
    $lpText = $lpText.Replace('Rule_0_3: x_C_3 <= 300', 'Rule_0_3: x_A_3 + x_B_3 >= 13')
    $lpText = $lpText.Replace('Rule_0_4: x_C_4 <= 300', 'Rule_0_4: x_A_4 + x_B_4 >= 13')
    $lpText = $lpText.Replace('Rule_0_5: x_C_5 <= 300', 'Rule_0_5: x_A_5 + x_B_5 >= 13')
    $lpText = $lpText.Replace('Rule_0_6: x_C_6 <= 300', 'Rule_0_6: x_A_6 + x_B_6 >= 13')
    $lpText = $lpText.Replace('Rule_0_7: x_C_7 <= 300', 'Rule_0_7: x_A_7 + x_B_7 >= 13')
    $lpText = $lpText.Replace('Rule_0_8: x_C_8 <= 300', 'Rule_0_8: x_A_8 + x_B_8 >= 13')
    $lpText = $lpText.Replace('Rule_0_9: x_C_9 <= 300', 'Rule_0_9: x_A_9 + x_B_9 >= 13')

$lpCell.Value = $lpText

Write-Output "done"
